# Commit: Cypress install & Read me doc update.
#
# 1) The "unit tests" bullet now calls out the actual testing
#    packages (Enzyme, Cypress, Jest) instead of the old "... js files."
#    wording.
# 2) The localStorage bullet drops the redundant "also" ("and also" -> "and").
$d = $word.ActiveDocument

$oldText1 = "Each component has their unit tests js files."
$newXml1 = '<w:p w14:paraId="6AA1EE99" w14:textId="78DCB2F4" w:rsidR="00AE02B2" w:rsidRDefault="00AE02B2" w:rsidP="00261FFA"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">For unit testing </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Enzyme, Cypress and Jes</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>t packages are used.</w:t></w:r></w:p>'

$oldText2 = "It makes use of localstorage to store the board values and also retrieve them on page load."
$newXml2 = '<w:p w14:paraId="4EE40B34" w14:textId="4E0C3D5B" w:rsidR="007D2D54" w:rsidRDefault="007D2D54" w:rsidP="00261FFA"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">It makes use of </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>localstorage</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> to store the board values </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>and</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> retrieve them on page load.</w:t></w:r></w:p>'

$updated = 0
foreach ($para in $d.Paragraphs) {
    $t = $para.Range.Text
    if ($t.StartsWith($oldText1)) {
        $para.Range.InsertXML($newXml1)
        $updated = $updated + 1
    }
}

foreach ($para in $d.Paragraphs) {
    $t = $para.Range.Text
    if ($t.StartsWith($oldText2)) {
        $para.Range.InsertXML($newXml2)
        $updated = $updated + 1
    }
}

if ($updated -ne 2) {
    throw "Expected to update 2 paragraphs, updated $updated"
}

Write-Output $d.Content.Text
